$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 283, shifting existing rows 283:391 down to 284:392.
$ws.Rows.Item(283).Insert()

# Populate the newly inserted row 283 with the new record's data.
$ws.Cells.Item(283, 1).Value = 10
$ws.Cells.Item(283, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(283, 3).Value = 'La Araucanía'
$ws.Cells.Item(283, 4).Value = 44524
$ws.Cells.Item(283, 5).Value = 9
$ws.Cells.Item(283, 6).Value = 100112027
$ws.Cells.Item(283, 7).Value = 'Melón'
$ws.Cells.Item(283, 8).Value = 'Tuna'
$ws.Cells.Item(283, 9).Value = 'Extra'
$ws.Cells.Item(283, 10).Value = 500
$ws.Cells.Item(283, 11).Value = 2000
$ws.Cells.Item(283, 12).Value = 2000
$ws.Cells.Item(283, 13).Value = 2000
$ws.Cells.Item(283, 14).Value = '$/unidad'
$ws.Cells.Item(283, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(283, 16).Value = 2000
$ws.Cells.Item(283, 17).Value = 1
$ws.Cells.Item(283, 18).Value = 'Hortaliza'
